$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the long citation strings in column G ("ref") with short citation keys.
# Rows 3, 6, 9, 10 previously referenced the Wong et al. citation.
$ws.Range("G3").Value = "wong2018assessment"
$ws.Range("G6").Value = "wong2018assessment"
$ws.Range("G9").Value = "wong2018assessment"
$ws.Range("G10").Value = "wong2018assessment"

# Rows 4 and 5 previously referenced the Latremouille-Viau et al. citation.
$ws.Range("G4").Value = "latremouille2017economic"
$ws.Range("G5").Value = "latremouille2017economic"

# The "notes" column (H) text is unchanged, but with the long citation strings
# gone the wrapped row heights need to be recomputed to fit the (now shorter)
# wrapped content in columns G/H.
$ws.Rows.Item(3).EntireRow.AutoFit()
$ws.Rows.Item(6).EntireRow.AutoFit()
$ws.Rows.Item(9).EntireRow.AutoFit()
$ws.Rows.Item(10).EntireRow.AutoFit()

# Rows 4 and 5 keep a taller, wrapped height (the long note in column H still
# wraps across several lines even though the G citation text got shorter).
$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 60
